$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and date range)
$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"

# Row 14: Murder
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 11
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = -45
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = -35.294117647058
$ws.Range("N14").Value = -91.2

# Row 15: Rape
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 74
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 48
$ws.Range("L15").Value = 48
$ws.Range("M15").Value = 4.225352112676
$ws.Range("N15").Value = -58.888888888888

# Row 16: Robbery
$ws.Range("C16").Value = 36
$ws.Range("D16").Value = 26
$ws.Range("E16").Value = 38.461538461538
$ws.Range("F16").Value = 130
$ws.Range("G16").Value = 153
$ws.Range("H16").Value = -15.032679738562
$ws.Range("I16").Value = 613
$ws.Range("J16").Value = 771
$ws.Range("K16").Value = -20.492866407263
$ws.Range("L16").Value = -7.121212121212
$ws.Range("M16").Value = -27.882352941176
$ws.Range("N16").Value = -82.247321169997

# Row 17: Fel. Assault
$ws.Range("C17").Value = 53
$ws.Range("D17").Value = 77
$ws.Range("E17").Value = -31.168831168831
$ws.Range("F17").Value = 243
$ws.Range("G17").Value = 286
$ws.Range("H17").Value = -15.034965034965
$ws.Range("I17").Value = 1048
$ws.Range("J17").Value = 1177
$ws.Range("K17").Value = -10.960067969413
$ws.Range("L17").Value = -1.964452759588
$ws.Range("M17").Value = 57.831325301204
$ws.Range("N17").Value = -47.756729810568

# Row 18: Burglary
$ws.Range("C18").Value = 37
$ws.Range("D18").Value = 25
$ws.Range("E18").Value = 48
$ws.Range("F18").Value = 109
$ws.Range("G18").Value = 100
$ws.Range("H18").Value = 9
$ws.Range("I18").Value = 528
$ws.Range("J18").Value = 531
$ws.Range("K18").Value = -0.564971751412
$ws.Range("L18").Value = -12.582781456953
$ws.Range("M18").Value = 4.761904761904
$ws.Range("N18").Value = -87.046123650637

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 118
$ws.Range("D19").Value = 139
$ws.Range("E19").Value = -15.107913669064
$ws.Range("F19").Value = 506
$ws.Range("G19").Value = 500
$ws.Range("H19").Value = 1.2
$ws.Range("I19").Value = 2170
$ws.Range("J19").Value = 2366
$ws.Range("K19").Value = -8.284023668639
$ws.Range("L19").Value = -1.898734177215
$ws.Range("M19").Value = 30.09592326139
$ws.Range("N19").Value = -45.871788475929

# Row 20: G.L.A.
$ws.Range("C20").Value = 25
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 82
$ws.Range("G20").Value = 61
$ws.Range("H20").Value = 34.426229508196
$ws.Range("I20").Value = 313
$ws.Range("J20").Value = 332
$ws.Range("K20").Value = -5.722891566265
$ws.Range("L20").Value = -35.1966873706
$ws.Range("M20").Value = 63.020833333333
$ws.Range("N20").Value = -91.39873591646

# Row 21: TOTAL
$ws.Range("C21").Value = 271
$ws.Range("D21").Value = 285
$ws.Range("E21").Value = -4.912280701754
$ws.Range("F21").Value = 1081
$ws.Range("G21").Value = 1114
$ws.Range("H21").Value = -2.962298025134
$ws.Range("I21").Value = 4757
$ws.Range("J21").Value = 5247
$ws.Range("K21").Value = -9.338669716028
$ws.Range("L21").Value = -6.725490196078
$ws.Range("M21").Value = 19.944528492183
$ws.Range("N21").Value = -72.798490393412

# Row 22: Transit
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -42.857142857142
$ws.Range("F22").Value = 15
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = -21.052631578947
$ws.Range("I22").Value = 84
$ws.Range("J22").Value = 102
$ws.Range("K22").Value = -17.647058823529
$ws.Range("L22").Value = -23.636363636363
$ws.Range("M22").Value = -4.545454545454
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("C23").Value = 26
$ws.Range("D23").Value = 24
$ws.Range("E23").Value = 8.333333333333
$ws.Range("F23").Value = 87
$ws.Range("G23").Value = 103
$ws.Range("H23").Value = -15.533980582524
$ws.Range("I23").Value = 479
$ws.Range("J23").Value = 490
$ws.Range("K23").Value = -2.244897959183
$ws.Range("L23").Value = 4.814004376367
$ws.Range("M23").Value = 64.041095890411
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("C24").Value = 269
$ws.Range("D24").Value = 246
$ws.Range("E24").Value = 9.349593495934
$ws.Range("F24").Value = 1044
$ws.Range("G24").Value = 959
$ws.Range("H24").Value = 8.863399374348
$ws.Range("I24").Value = 5293
$ws.Range("J24").Value = 4663
$ws.Range("K24").Value = 13.510615483594
$ws.Range("L24").Value = 3.96778628953
$ws.Range("M24").Value = 63.768564356435
$ws.Range("N24").Value = "***.*"

# Row 25: Retail Theft
$ws.Range("C25").Value = 120
$ws.Range("D25").Value = 124
$ws.Range("E25").Value = -3.225806451612
$ws.Range("F25").Value = 513
$ws.Range("G25").Value = 496
$ws.Range("H25").Value = 3.427419354838
$ws.Range("I25").Value = 2928
$ws.Range("J25").Value = 2486
$ws.Range("K25").Value = 17.779565567176
$ws.Range("L25").Value = 2.485124256212
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"

# Row 26: Misd. Assault
$ws.Range("C26").Value = 90
$ws.Range("D26").Value = 118
$ws.Range("E26").Value = -23.728813559322
$ws.Range("F26").Value = 439
$ws.Range("G26").Value = 413
$ws.Range("H26").Value = 6.295399515738
$ws.Range("I26").Value = 1776
$ws.Range("J26").Value = 1807
$ws.Range("K26").Value = -1.715550636413
$ws.Range("L26").Value = 7.117008443908
$ws.Range("M26").Value = -6.820566631689
$ws.Range("N26").Value = "***.*"

# Row 27: UCR Rape*
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = -27.777777777777
$ws.Range("I27").Value = 92
$ws.Range("J27").Value = 79
$ws.Range("K27").Value = 16.455696202531
$ws.Range("L27").Value = -2.127659574468
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28: Other Sex Crimes
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 14
$ws.Range("E28").Value = 57.142857142857
$ws.Range("F28").Value = 69
$ws.Range("G28").Value = 50
$ws.Range("H28").Value = 38
$ws.Range("I28").Value = 230
$ws.Range("J28").Value = 215
$ws.Range("K28").Value = 6.976744186046
$ws.Range("L28").Value = 4.072398190045
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"

# Row 29: Shooting Vic.
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -25
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = 12
$ws.Range("H29").Value = -16.666666666666
$ws.Range("I29").Value = 30
$ws.Range("J29").Value = 41
$ws.Range("K29").Value = -26.829268292682
$ws.Range("L29").Value = -48.275862068965
$ws.Range("M29").Value = -55.223880597014
$ws.Range("N29").Value = -89.690721649484

# Row 30: Shooting Inc.
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = -10
$ws.Range("I30").Value = 27
$ws.Range("J30").Value = 33
$ws.Range("K30").Value = -18.181818181818
$ws.Range("L30").Value = -49.056603773584
$ws.Range("M30").Value = -55.737704918032
$ws.Range("N30").Value = -89.96282527881

# Row 31: Hate Crimes
$ws.Range("C31").Value = "0"
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = -87.5
$ws.Range("I31").Value = 25
$ws.Range("J31").Value = 44
$ws.Range("K31").Value = -43.181818181818
$ws.Range("L31").Value = -10.714285714285
$ws.Range("M31").Value = "***.*"
$ws.Range("N31").Value = "***.*"

# Row 33: Traffic Fatalities
$ws.Range("C33").Value = "0"
$ws.Range("D33").Value = "0"
$ws.Range("E33").Value = "***.*"
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = "0"
$ws.Range("H33").Value = "***.*"
$ws.Range("I33").Value = 7
$ws.Range("J33").Value = 5
$ws.Range("K33").Value = 40
$ws.Range("L33").Value = -30
$ws.Range("M33").Value = "***.*"
$ws.Range("N33").Value = "***.*"

